$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.013269782066345
$ws.Range("B1").Value = 1.579648494720459
$ws.Range("C1").Value = 3.725919723510742
$ws.Range("D1").Value = 3.026717662811279
$ws.Range("E1").Value = 1.421087861061096
